# Applies the diff: normalize the subject metadata for sub-1 and sub-2,
# drop the duplicate sub-2 row (formerly row 4, originally listed as sub-3
# in row 3), so the sheet ends with just two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (currently "sub-3") becomes "sub-2" with its own updated values,
# effectively replacing/merging with the old row 4 data before that row
# is removed entirely.
$ws.Range("A3").Value = "sub-2"

# Row 2 (sub-1) field updates
$ws.Range("C2").Value = "Control"
$ws.Range("D2").Value = "30y"
$ws.Range("F2").Value = "Human"
$ws.Range("G2").Value = "Not Defined"
$ws.Range("H2").Value = "Not Defined"
$ws.Range("I2").Value = "Prime Adult Stage"
$ws.Range("J2").Value = "Not Defined"
$ws.Range("K2").Value = "Not Defined"

# Row 3 (sub-2) field updates
$ws.Range("C3").Value = "Control"
$ws.Range("D3").Value = "20y"
$ws.Range("F3").Value = "Human"
$ws.Range("G3").Value = "Not Defined"
$ws.Range("H3").Value = "Not Defined"
$ws.Range("I3").Value = "Prime Adult Stage"
$ws.Range("J3").Value = "Not Defined"
$ws.Range("K3").Value = "Not Defined"

# Remove the old row 4 (its data has been folded into row 3 above).
$ws.Rows.Item(4).Delete()
